$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.468.50"
$ws.Range("E2").Value = "  -0.32%  "
$ws.Range("D3").Value = "3.508.43"
$ws.Range("E3").Value = "  -3.11%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "604.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.23"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.06%  "
$ws.Range("D7").Value = "3.509.58"
$ws.Range("E7").Value = "  -3.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.504"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.70"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.18%  "
$ws.Range("E11").Value = "  -5.89%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.402"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.50%  "
$ws.Range("D13").Value = "4.095.57"
$ws.Range("E13").Value = "  -3.27%  "
$ws.Range("E14").Value = "  -7.88%  "
$ws.Range("E15").Value = "  -5.07%  "
$ws.Range("D16").Value = "3.490.02"
$ws.Range("E16").Value = "  -3.45%  "
$ws.Range("E17").Value = "  -0.40%  "
$ws.Range("D18").Value = "66.334.13"
$ws.Range("E19").Value = "  -10.21%  "
$ws.Range("E20").Value = "  -4.89%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.54"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.84%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "419.84"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.98%  "
$ws.Range("E23").Value = "  -5.68%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "76.66"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.87%  "
$ws.Range("D25").Value = "3.651.20"
$ws.Range("E25").Value = "  -3.01%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("E27").Value = "  -9.28%  "
$ws.Range("E28").Value = "  -3.41%  "
$ws.Range("E29").Value = "  -7.73%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.72"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -7.99%  "
$ws.Range("E31").Value = "  +0.13%  "
$ws.Range("D32").Value = "3.513.02"
$ws.Range("E32").Value = "  -2.90%  "
$ws.Range("E33").Value = "  -4.36%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "24.09"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.38%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.32"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -10.74%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "7.44"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.82%  "
$ws.Range("E38").Value = "  -4.99%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "173.62"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.04%  "
$ws.Range("E40").Value = "  -8.57%  "
$ws.Range("E41").Value = "  -7.25%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.92"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.849"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.63%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "45.45"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.76%  "
$ws.Range("E45").Value = "  -7.93%  "
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("E47").Value = "  -10.94%  "
$ws.Range("E48").Value = "  -3.03%  "
$ws.Range("B49").Value = "ONDO"
$ws.Range("C49").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.10"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.37%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.77"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.13%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.891"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -8.09%  "
